$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$ws.Range("D2").Value = 3445.9960718981852
$ws.Range("E2").Value = 3440.5901179619041
$ws.Range("F2").Value = 3074.847792316134
$ws.Range("G2").Value = 1572.1875028153099
$ws.Range("H2").Value = 2978.431393750288
$ws.Range("I2").Value = 2395.7918443106259
$ws.Range("J2").Value = 1822.1332839254119
$ws.Range("K2").Value = 991.47366189812465
$ws.Range("L2").Value = 618.81038201592662
$ws.Range("D3").Value = 181.3682143104308
$ws.Range("E3").Value = 10.159250988713261
$ws.Range("F3").Value = 147.3939683815538
$ws.Range("G3").Value = 808.32262355542946
$ws.Range("H3").Value = 277.82070714077929
$ws.Range("K3").Value = 139.91327045606249
$ws.Range("L3").Value = 518.21838555889519
$ws.Range("F4").Value = 12.77807870661864
$ws.Range("G4").Value = 44.457744295548629
$ws.Range("H4").Value = 89.21238632578266
$ws.Range("I4").Value = 86.894004715929356
$ws.Range("J4").Value = 82.535653975715562
$ws.Range("K4").Value = 61.639709847064523
$ws.Range("L4").Value = 72.576304313286499
$ws.Range("D5").Value = 128.41592936589771
$ws.Range("E5").Value = 46.742200733010243
$ws.Range("F5").Value = 39.647569577502942
$ws.Range("G5").Value = 858.4831749185696
$ws.Range("H5").Value = 282.29645907281991
$ws.Range("I5").Value = 822.61854482054957
$ws.Range("J5").Value = 1.976826965537904
$ws.Range("K5").Value = 2.0156901200618251
$ws.Range("E6").Value = 932.39641152252648
$ws.Range("F6").Value = 824.43887914529375
$ws.Range("G6").Value = 679.56576076054205
$ws.Range("H6").Value = 71.397636847318438
$ws.Range("L6").Value = 711.55658708203805
$ws.Range("D7").Value = 71.397636847318438
$ws.Range("G7").Value = 171.23986492528331
$ws.Range("H7").Value = 868.38984346784946
$ws.Range("I7").Value = 1001.949075647482
$ws.Range("J7").Value = 1005.302283046533
$ws.Range("K7").Value = 983.2467772242087
$ws.Range("L7").Value = 278.79600825025562
$ws.Range("D8").Value = 305.70915085901407
$ws.Range("E8").Value = 345.83483815730108
$ws.Range("F8").Value = 366.23381126648928
$ws.Range("G8").Value = 397.26306891645572
$ws.Range("H8").Value = 410.8659830687098
$ws.Range("I8").Value = 420.26022820029129
$ws.Range("J8").Value = 439.46758315384761
$ws.Range("K8").Value = 392.63590467200078
$ws.Range("L8").Value = 391.66851947575759
$ws.Range("D9").Value = 41.058446613978532
$ws.Range("E9").Value = 33.023917785009552
$ws.Range("F9").Value = 24.98938895604033
$ws.Range("G9").Value = 16.95486012707633
$ws.Range("H9").Value = 8.9324477834599207
$ws.Range("I9").Value = 102.9254707395301
$ws.Range("J9").Value = 103.2874842322735
$ws.Range("K9").Value = 92.55724829562763
$ws.Range("D10").Value = 90.326119315574786
$ws.Range("E10").Value = 70.530965414970851
$ws.Range("F10").Value = 54.500805198475753
$ws.Range("G10").Value = 38.586858613015167
$ws.Range("H10").Value = 23.45410118689351
$ws.Range("I10").Value = 8.9224299224256747
$ws.Range("J10").Value = 7.0398985378976926
$ws.Range("K10").Value = 6.7089628334579459
$ws.Range("L10").Value = 6.6936467118908594
$ws.Range("D11").Value = 1.1576670811117009
$ws.Range("E11").Value = 2.9387902256237859
$ws.Range("F11").Value = 2.2880218414381308
$ws.Range("G11").Value = 1.6173761607122821
$ws.Range("H11").Value = 1.0016590268338459
$ws.Range("I11").Value = 0.37868561811768248
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("D12").Value = 419.46128716965188
$ws.Range("D13").Value = 36.474894536491469
$ws.Range("E13").Value = 267.77839246018999
$ws.Range("F13").Value = 1088.6615017154149
$ws.Range("G13").Value = 869.28157828225267
$ws.Range("H13").Value = 218.9450244384754
$ws.Range("I13").Value = 150.7356559359649
$ws.Range("L13").Value = 137.25427606500239
$ws.Range("D14").Value = 643.15511869071145
$ws.Range("E14").Value = 554.88438158701797
$ws.Range("F14").Value = 447.40711373049862
$ws.Range("G14").Value = 313.05510834820842
$ws.Range("H14").Value = 347.27504552087208
$ws.Range("I14").Value = 297.53435126778021
$ws.Range("J14").Value = 251.98347954532801
$ws.Range("D15").Value = 83.573828982408841
$ws.Range("E15").Value = 72.103620206222672
$ws.Range("F15").Value = 58.137647546900929
$ws.Range("G15").Value = 40.679477355981881
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("I16").Value = 38.662655814457302
$ws.Range("J16").Value = 32.743615986116083
$ws.Range("K16").Value = 33.176519490408467
$ws.Range("E17").Value = 38.437859370490671
$ws.Range("F17").Value = 40.335944287880316
$ws.Range("G17").Value = 293.22796655669111
$ws.Range("H17").Value = 92.179999999999978
$ws.Range("I17").Value = 92.179999999999978
$ws.Range("J17").Value = 70.585255168069054
$ws.Range("L17").Value = 0
$ws.Range("D18").Value = 41.538924622552173
$ws.Range("E18").Value = 48.92091192607905
$ws.Range("F18").Value = 51.336656366393157
$ws.Range("G18").Value = 373.19923016306149
$ws.Range("H18").Value = 448.33852220409602
$ws.Range("I18").Value = 444.39689689557218
$ws.Range("J18").Value = 89.83577930481519
$ws.Range("K18").Value = 60.279352327418138
$ws.Range("L18").Value = 36.046367523626976
$ws.Range("D19").Value = 345.80307436793771
$ws.Range("E19").Value = 307.75617322383277
$ws.Range("F19").Value = 115.6528970119132
$ws.Range("G19").Value = 169.43096843470781
$ws.Range("H19").Value = 10.827934154582691
$ws.Range("I19").Value = 0.48325471388551472
$ws.Range("J19").Value = 0.49420674138454962
$ws.Range("K19").Value = 0.50392253001544596
$ws.Range("L19").Value = 0.51240730679663216
$ws.Range("D20").Value = 183.85117144902509
$ws.Range("E20").Value = 147.37435279708029
$ws.Range("F20").Value = 110.8975341451351
$ws.Range("G20").Value = 94.803232809996743
$ws.Range("H20").Value = 37.943896841245227
$ws.Range("D21").Value = 160.31932720355559
$ws.Range("E21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("G22").Value = 50.391892405358178
$ws.Range("H22").Value = 63.828455708062457
$ws.Range("I22").Value = 36.26435803984954
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 142.32087760822611
$ws.Range("F23").Value = 122.4013839495879
$ws.Range("G23").Value = 784.53911155837363
$ws.Range("H23").Value = 1073.0004544340391
$ws.Range("I23").Value = 1590.833987631587
$ws.Range("J23").Value = 1136.739932323502
$ws.Range("K23").Value = 1590.1511693481009
$ws.Range("L23").Value = 1422.663652594583
$ws.Range("E24").Value = 70.06236026914732
$ws.Range("F24").Value = 74.441782657444904
$ws.Range("G24").Value = 80.259679796609575
$ws.Range("H24").Value = 82.918789392896599
$ws.Range("I24").Value = 85.074077905040056
$ws.Range("J24").Value = 82.774317489707727
$ws.Range("K24").Value = 96.21480830755867
$ws.Range("L24").Value = 82.838303471857529
$ws.Range("D25").Value = 4482.0940540540541
$ws.Range("E25").Value = 4481.2205405405412
$ws.Range("G25").Value = 4404.6833688119241
$ws.Range("H25").Value = 4607.8654269770304
$ws.Range("I25").Value = 4518.5533054450871
$ws.Range("J25").Value = 6348.6671974504516
$ws.Range("K25").Value = 7252.6797792795978
$ws.Range("L25").Value = 7648.7004609290088
$ws.Range("D26").Value = 1638.1351596571919
$ws.Range("E26").Value = 1517.414143703676
$ws.Range("F26").Value = 915.71327362566456
$ws.Range("G26").Value = 963.34089147986469
$ws.Range("H26").Value = 892.73063460857202
$ws.Range("I26").Value = 850.30516129451837
$ws.Range("J26").Value = 826.03839020791384
$ws.Range("K26").Value = 1046.490362955893
$ws.Range("L26").Value = 1138.502652541416
$ws.Range("D27").Value = 1402.830166589318
$ws.Range("E27").Value = 1012.425935375055
$ws.Range("F27").Value = 981.7648986214939
$ws.Range("G27").Value = 825.29046187451604
$ws.Range("H27").Value = 768.8701196399735
$ws.Range("I27").Value = 617.04381523413895
$ws.Range("J27").Value = 1500.614443140603
$ws.Range("K27").Value = 1745.7846385950229
$ws.Range("L27").Value = 1835.9681794583689
$ws.Range("D28").Value = 0
$ws.Range("H28").Value = 128.23099999999999
$ws.Range("I28").Value = 115.36607968818861
$ws.Range("J28").Value = 129.2934180234393
$ws.Range("K28").Value = 128.85492894907031
$ws.Range("E29").Value = 40.10079908079102
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 1.21954609062024
$ws.Range("J29").Value = 46.121280028259527
$ws.Range("D30").Value = 128.23099999999999
$ws.Range("E30").Value = 128.23099999999999
$ws.Range("F30").Value = 128.23099999999999
$ws.Range("G30").Value = 128.23099999999999
$ws.Range("E31").Value = 83.905363152681446
$ws.Range("F31").Value = 148.33422532479199
$ws.Range("G31").Value = 207.34600772868981
$ws.Range("H31").Value = 147.06927950629299
$ws.Range("D32").Value = 82.905363152681446
$ws.Range("E32").Value = 147.33422532479199
$ws.Range("F32").Value = 29.355895529914282
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 84.966731056309996
$ws.Range("D33").Value = 1828.3341973495019
$ws.Range("E33").Value = 2128.069459436967
$ws.Range("F33").Value = 1310.5962437258311
$ws.Range("G33").Value = 57.008658444249598
$ws.Range("H33").Value = 28.338930934761169
$ws.Range("I33").Value = 12.186463418434901
$ws.Range("L33").Value = -0.00000001108756737266958
$ws.Range("D34").Value = 168.2857328926988
$ws.Range("E34").Value = 359.69721637893838
$ws.Range("F34").Value = 372.11081639724631
$ws.Range("G34").Value = 1544.5546619874731
$ws.Range("H34").Value = 2168.775345499133
$ws.Range("I34").Value = 1017.934293860725
$ws.Range("J34").Value = 669.00402330540317
$ws.Range("K34").Value = 251.3830518659515
$ws.Range("L34").Value = 150.32420765823409
$ws.Range("G35").Value = 289.44833784240882
$ws.Range("H35").Value = 931.3639355457301
$ws.Range("I35").Value = 2606.151709913303
$ws.Range("J35").Value = 2639.515032420737
$ws.Range("K35").Value = 2675.960886434907
$ws.Range("L35").Value = 2702.6168919477768
$ws.Range("D36").Value = 397.082771650609
$ws.Range("E36").Value = 389.83686144161692
$ws.Range("F36").Value = 627.789461217802
$ws.Range("G36").Value = 561.57734454679905
$ws.Range("H36").Value = 132.37434788641011
$ws.Range("I36").Value = 83.318220322157558
$ws.Range("J36").Value = 366.45395109275609
$ws.Range("K36").Value = 40.399012414310413
$ws.Range("L36").Value = 1.008064516129032
$ws.Range("D37").Value = 20.899093244768899
$ws.Range("E37").Value = 55.690980205945252
$ws.Range("F37").Value = 158.93404081463339
$ws.Range("G37").Value = 24.984274819269821
$ws.Range("H37").Value = 18.752184520911459
$ws.Range("K37").Value = 428.72421337635518
$ws.Range("L37").Value = 480.11082833911041
$ws.Range("F38").Value = 7.9467020407316724
$ws.Range("G38").Value = 10.954461482218599
$ws.Range("H38").Value = 4.140452942666343
$ws.Range("I38").Value = 3.021904364016077
$ws.Range("J38").Value = 16.598959457163161
$ws.Range("K38").Value = 25.55855887435963
$ws.Range("L38").Value = 30.709716565228039
$ws.Range("F39").Value = 10.125589393703761
$ws.Range("D40").Value = 194.28994483416469
$ws.Range("E40").Value = 160.5699688987703
$ws.Range("F40").Value = 154.00888240369611
$ws.Range("G40").Value = 194.28994483416469
$ws.Range("H40").Value = 141.22228065382961
$ws.Range("I40").Value = 23.018875345932599
$ws.Range("D41").Value = 1310.7659256320619
$ws.Range("E41").Value = 1109.356826776168
$ws.Range("F41").Value = 1303.167102988087
$ws.Range("G41").Value = 1090.3479152570351
$ws.Range("H41").Value = 305.86389511470833
$ws.Range("I41").Value = 131.68150749364389
$ws.Range("K41").Value = 0
$ws.Range("D42").Value = 241.73283834018159
$ws.Range("E42").Value = 199.7566874990645
$ws.Range("F42").Value = 191.57202445337811
$ws.Range("G42").Value = 132.86695395474811
$ws.Range("H42").Value = 66.449907019454656
$ws.Range("I42").Value = 28.773594182415739
